$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + data values in column O ("My Appointments"), copying style from
# the adjacent column N so the new cells match the existing header/data formatting.
$ws.Range("N2").Copy($ws.Range("O2"))
$ws.Range("O2").Value = "My Appointments"

$ws.Range("N3").Copy($ws.Range("O3"))
$ws.Range("O3").Value = "Yes"

$ws.Range("N4").Copy($ws.Range("O4"))
$ws.Range("O4").Value = "No"

# Existing Yes/No values that flipped.
$ws.Range("B3").Value = "Yes"
$ws.Range("I3").Value = "No"
$ws.Range("I4").Value = "No"

# Column width adjustments triggered by the new column / wider header text.
$ws.Columns("N").ColumnWidth = 8
$ws.Columns("O").ColumnWidth = 15.5

# Selection moved.
[void]$ws.Range("E13").Select()
